# Implemented passing authentication via the request headers
# Adds a new "Create Hospital" API row (with a hyperlink on the endpoint
# cell) to the API_Doc worksheet, and makes column B wrap-text / wider to
# match the rest of the documentation table's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wrap text on the existing header/body cells in column B (B1, B2) ---
# (these reuse the workbook's existing "wrap text" cell style, the same
# one already used by C2)
$ws.Range("B1").WrapText = $true
$ws.Range("B2").WrapText = $true

# --- New row 3: Create Hospital endpoint ---
$ws.Range("A3").Value = "Create Hospital"

# B3 becomes a hyperlink pointing at the create-hospital endpoint; Add()
# also writes the display text (the URL) into the cell. Wrap text is
# turned on first so the generated Hyperlink-style xf picks up the
# existing wrap-text cell style (xfId=1) instead of a fresh one.
$ws.Range("B3").WrapText = $true
[void]$ws.Hyperlinks.Add($ws.Range("B3"), "http://ramsayportalapi-uat.azurewebsites.net/api/hospitals/create")

$payload = @"
{
 "HospitalId":#hospitalId,
 "HospitalName":"#hospitalName",
 "Address1":"dddd",
 "Address2":"ddddd",
 "Suburb":"adf",
 "Web":"tyntymtwym",
 "SiteCode":"ddd",
 "Postcode":1,
 "State":"ddd",
 "AdminEmail":"sdcs@hadvc.com",
 "Longitude":0,
 "Latitude":0,
 "IsActive":true,
 "CreateUser":"portaladmin"
}
"@
$ws.Range("C3").Value = $payload
$ws.Range("C3").WrapText = $true

# Row 3 is tall enough to show the whole wrapped JSON payload.
$ws.Rows.Item(3).RowHeight = 240

# --- Column sizing to fit the new, wider content ---
# ColumnWidth is specified in characters; the stored (pixel-grid) width it
# produces is what ends up in the saved file.
$ws.Columns.Item(1).ColumnWidth = 13.833333333333332
$ws.Columns.Item(2).ColumnWidth = 66.83333333333334
$ws.Columns.Item(3).ColumnWidth = 34.166666666666664

# --- Final selection, matching where the author's cursor ended up ---
[void]$ws.Range("C3").Select()
